# Applies the "valid and invalid with excel and dataprovider" edit:
#  - invalidCredentialTest: append a new data row (Har / admin123 / Invalid credentials)
#  - validCredentialTest: remove the last data row (Balaji / admin123 / hyperlink) and
#    its hyperlink
#  - switch the active sheet/tab from invalidCredentialTest to validCredentialTest and
#    update each sheet's saved selection accordingly

$wb = $excel.ActiveWorkbook

$wsInvalid = $wb.Worksheets.Item("invalidCredentialTest")
$wsValid   = $wb.Worksheets.Item("validCredentialTest")

# --- invalidCredentialTest: add new row 4 ---
$wsInvalid.Range("A4").Value = "Har"
$wsInvalid.Range("B4").Value = "admin123"
$wsInvalid.Range("C4").Value = "Invalid credentials"
$wsInvalid.Range("C4").Select()

# --- validCredentialTest: delete row 3 (data + hyperlink) ---
# (iterate so the hyperlink that actually sits on C3 gets removed, leaving
#  the C2 hyperlink/relationship untouched)
$i = 0
foreach ($h in $wsValid.Hyperlinks) {
    $i++
    if ($i -eq 2) {
        $h.Delete()
    }
}
$wsValid.Rows.Item(3).Delete()

# select the entire column D on the valid-credentials sheet
$wsValid.Range("D1:XFD1048576").Select()

# make validCredentialTest the active/visible tab
$wsValid.Activate()
